# Add working set of sequences
# Rewrites the data rows (image/word/category/random-number columns)
# of the cue sequence sheet with a new working set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 82
$ws.Range("C2").Value = "house/house030.jpg"
$ws.Range("D2").Value = "laufen"
$ws.Range("E2").Value = "house"
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = "face/face017.jpg"
$ws.Range("D3").Value = "währen"
$ws.Range("E3").Value = "face"
$ws.Range("B4").Value = 112
$ws.Range("C4").Value = "face/face009.jpg"
$ws.Range("D4").Value = "stechen"
$ws.Range("E4").Value = "face"
$ws.Range("B5").Value = 120
$ws.Range("C5").Value = "face/face007.jpg"
$ws.Range("D5").Value = "mieten"
$ws.Range("E5").Value = "face"
$ws.Range("B6").Value = 39
$ws.Range("C6").Value = "face/face015.jpg"
$ws.Range("D6").Value = "saufen"
$ws.Range("E6").Value = "face"
$ws.Range("B7").Value = 94
$ws.Range("C7").Value = "house/house024.jpg"
$ws.Range("D7").Value = "hupen"
$ws.Range("E7").Value = "house"
$ws.Range("B8").Value = 83
$ws.Range("C8").Value = "face/face026.jpg"
$ws.Range("D8").Value = "dauern"
$ws.Range("E8").Value = "face"
$ws.Range("B9").Value = 109
$ws.Range("C9").Value = "face/face008.jpg"
$ws.Range("D9").Value = "schmecken"
$ws.Range("E9").Value = "face"
$ws.Range("B10").Value = 65
$ws.Range("C10").Value = "face/face028.jpg"
$ws.Range("D10").Value = "fesseln"
$ws.Range("E10").Value = "face"
$ws.Range("B11").Value = 127
$ws.Range("C11").Value = "face/face002.jpg"
$ws.Range("D11").Value = "füllen"
$ws.Range("E11").Value = "face"
$ws.Range("B12").Value = 37
$ws.Range("C12").Value = "house/house015.jpg"
$ws.Range("D12").Value = "jubeln"
$ws.Range("E12").Value = "house"
$ws.Range("B13").Value = 126
$ws.Range("C13").Value = "house/house013.jpg"
$ws.Range("D13").Value = "opfern"
$ws.Range("E13").Value = "house"
$ws.Range("B14").Value = 9
$ws.Range("C14").Value = "house/house000.jpg"
$ws.Range("D14").Value = "strahlen"
$ws.Range("E14").Value = "house"
$ws.Range("B15").Value = 110
$ws.Range("C15").Value = "face/face000.jpg"
$ws.Range("D15").Value = "drohen"
$ws.Range("E15").Value = "face"
$ws.Range("B16").Value = 43
$ws.Range("C16").Value = "house/house029.jpg"
$ws.Range("D16").Value = "kehren"
$ws.Range("E16").Value = "house"
$ws.Range("B17").Value = 33
$ws.Range("C17").Value = "house/house004.jpg"
$ws.Range("D17").Value = "stärken"
$ws.Range("E17").Value = "house"
$ws.Range("B18").Value = 72
$ws.Range("C18").Value = "house/house001.jpg"
$ws.Range("D18").Value = "bleiben"
$ws.Range("E18").Value = "house"
$ws.Range("B19").Value = 48
$ws.Range("C19").Value = "house/house005.jpg"
$ws.Range("D19").Value = "starten"
$ws.Range("E19").Value = "house"
$ws.Range("B20").Value = 29
$ws.Range("C20").Value = "house/house007.jpg"
$ws.Range("D20").Value = "segeln"
$ws.Range("E20").Value = "house"
$ws.Range("B21").Value = 16
$ws.Range("C21").Value = "face/face031.jpg"
$ws.Range("D21").Value = "fliehen"
$ws.Range("E21").Value = "face"
$ws.Range("B22").Value = 24
$ws.Range("C22").Value = "face/face016.jpg"
$ws.Range("D22").Value = "sondern"
$ws.Range("E22").Value = "face"
$ws.Range("B23").Value = 108
$ws.Range("C23").Value = "house/house019.jpg"
$ws.Range("D23").Value = "bitten"
$ws.Range("E23").Value = "house"
$ws.Range("B24").Value = 27
$ws.Range("C24").Value = "house/house031.jpg"
$ws.Range("D24").Value = "backen"
$ws.Range("E24").Value = "house"
$ws.Range("B25").Value = 111
$ws.Range("C25").Value = "house/house027.jpg"
$ws.Range("D25").Value = "raten"
$ws.Range("E25").Value = "house"
$ws.Range("B26").Value = 60
$ws.Range("C26").Value = "face/face012.jpg"
$ws.Range("D26").Value = "füttern"
$ws.Range("E26").Value = "face"
$ws.Range("B27").Value = 51
$ws.Range("C27").Value = "house/house028.jpg"
$ws.Range("D27").Value = "schenken"
$ws.Range("E27").Value = "house"
$ws.Range("B28").Value = 124
$ws.Range("C28").Value = "face/face001.jpg"
$ws.Range("D28").Value = "töten"
$ws.Range("E28").Value = "face"
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = "face/face030.jpg"
$ws.Range("D29").Value = "scheitern"
$ws.Range("E29").Value = "face"
$ws.Range("B30").Value = 95
$ws.Range("C30").Value = "face/face023.jpg"
$ws.Range("D30").Value = "gründen"
$ws.Range("E30").Value = "face"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = "face/face024.jpg"
$ws.Range("D31").Value = "hauen"
$ws.Range("E31").Value = "face"
$ws.Range("B32").Value = 34
$ws.Range("C32").Value = "house/house018.jpg"
$ws.Range("D32").Value = "ehren"
$ws.Range("E32").Value = "house"
$ws.Range("B33").Value = 7
$ws.Range("C33").Value = "house/house017.jpg"
$ws.Range("D33").Value = "lehnen"
$ws.Range("E33").Value = "house"
